$d = $word.ActiveDocument

$replacements = @(
    @("491÷3=", "727÷6="),
    @("888÷6=", "896÷5="),
    @("249÷3=", "428÷5="),
    @("502÷3=", "228÷3="),
    @("172÷8=", "747÷6="),
    @("313÷8=", "444÷8="),
    @("684÷7=", "156÷7="),
    @("121÷9=", "306÷4="),
    @("829÷2=", "989÷6="),
    @("157÷7=", "646÷7="),
    @("830÷4=", "510÷7="),
    @("743÷2=", "488÷7="),
    @("979÷2=", "782÷7="),
    @("762÷5=", "771÷4="),
    @("748÷9=", "194÷5="),
    @("808÷8=", "586÷5="),
    @("692÷7=", "632÷7="),
    @("462÷3=", "813÷9="),
    @("725÷5=", "942÷3="),
    @("482÷4=", "506÷2="),
    @("555÷3=", "411÷5="),
    @("797÷4=", "980÷9="),
    @("359÷8=", "802÷6="),
    @("815÷3=", "304÷2="),
    @("370÷5=", "892÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
